# Generate Report for Archive
#
# 1) Update status text from "Ready for handoff" to "In Translation" on all
#    sheets/cells that currently show it (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2) Narrow the "Status" columns (Overview E:F, zh-cn C, de-de C) from their
#    previous wider width down to match the other report columns.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Compare against $cell.Text (always a string) and put the known
        # string literal on the left of -eq so PowerShell never coerces it
        # to a boolean when the cell happens to hold "True"/"False".
        if ($oldStatus -eq $cell.Text) {
            $cell.Value2 = $newStatus
        }
    }
}

# Narrow the previously wide "Status" columns. The new stored width matches
# what the other (already-narrow) columns use, set via ColumnWidth so it
# round-trips the same way Excel itself persists column widths.
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
